$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: columns A (Case) and B (End Time) for rows 32..41
$data = @(
    @(30, 100.8333333333333),
    @(31, 83.93333333333334),
    @(32, 149.7),
    @(33, 91.90000000000001),
    @(34, 101.4833333333333),
    @(35, 42.81666666666667),
    @(36, 45.48333333333333),
    @(37, 37.53333333333333),
    @(38, 42.76666666666667),
    @(39, 117.6166666666667)
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
